$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on Price cells whose new values would otherwise
# be auto-converted to numbers (losing literal formatting like trailing zeros).
$ws.Range("D4:D9").NumberFormat = "@"
$ws.Range("D11:D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19:D23").NumberFormat = "@"
$ws.Range("D25:D29").NumberFormat = "@"
$ws.Range("D31:D48").NumberFormat = "@"
$ws.Range("D50:D51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) values
$ws.Range("D2").Value = "23.766.72"
$ws.Range("E2").Value = "  +15.57%  "
$ws.Range("D3").Value = "1.653.95"
$ws.Range("E3").Value = "  +12.50%  "
$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  -1.20%  "
$ws.Range("D5").Value = "305.69"
$ws.Range("E5").Value = "  +10.87%  "
$ws.Range("D6").Value = "0.9847"
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("D7").Value = "0.3707"
$ws.Range("E7").Value = "  +4.36%  "
$ws.Range("D8").Value = "0.3432"
$ws.Range("E8").Value = "  +12.05%  "
$ws.Range("D9").Value = "47.42"
$ws.Range("E9").Value = "  +20.86%  "
$ws.Range("E10").Value = "  +7.78%  "
$ws.Range("D11").Value = "0.07171"
$ws.Range("E11").Value = "  +8.49%  "
$ws.Range("D12").Value = "0.9868"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "20.50"
$ws.Range("E13").Value = "  +13.77%  "
$ws.Range("D14").Value = "5.974"
$ws.Range("E14").Value = "  +9.76%  "
$ws.Range("D15").Value = "6.704"
$ws.Range("E15").Value = "  +8.78%  "
$ws.Range("D16").Value = "1.660.28"
$ws.Range("E16").Value = "  +12.78%  "
$ws.Range("D17").Value = "0.00001091"
$ws.Range("E17").Value = "  +7.39%  "
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").Value = "0.06713"
$ws.Range("E19").Value = "  +13.24%  "
$ws.Range("D20").Value = "80.74"
$ws.Range("E20").Value = "  +17.87%  "
$ws.Range("D21").Value = "16.33"
$ws.Range("E21").Value = "  +13.36%  "
$ws.Range("D22").Value = "6.054"
$ws.Range("E22").Value = "  +10.87%  "
$ws.Range("D23").Value = "11.88"
$ws.Range("E23").Value = "  +6.57%  "
$ws.Range("D24").Value = "23.852.74"
$ws.Range("E24").Value = "  +15.98%  "
$ws.Range("D25").Value = "2.327"
$ws.Range("E25").Value = "  +2.69%  "
$ws.Range("D26").Value = "3.408"
$ws.Range("E26").Value = "  -8.42%  "
$ws.Range("D27").Value = "2.656"
$ws.Range("E27").Value = "  +27.99%  "
$ws.Range("D28").Value = "152.05"
$ws.Range("E28").Value = "  +4.35%  "
$ws.Range("D29").Value = "19.40"
$ws.Range("E29").Value = "  +13.65%  "
$ws.Range("D30").Value = "1.845.64"
$ws.Range("E30").Value = "  +13.12%  "
$ws.Range("D31").Value = "125.90"
$ws.Range("E31").Value = "  +10.25%  "
$ws.Range("D32").Value = "4.024"
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("D33").Value = "6.100"
$ws.Range("E33").Value = "  +24.65%  "
$ws.Range("D34").Value = "0.9795"
$ws.Range("E34").Value = "  +24.35%  "
$ws.Range("D35").Value = "1.703"
$ws.Range("E35").Value = "  +17.43%  "
$ws.Range("D36").Value = "0.08347"
$ws.Range("E36").Value = "  +5.38%  "
$ws.Range("D37").Value = "12.16"
$ws.Range("E37").Value = "  +18.91%  "
$ws.Range("D38").Value = "8.852"
$ws.Range("E38").Value = "  +21.64%  "
$ws.Range("D39").Value = "0.06298"
$ws.Range("E39").Value = "  +11.18%  "
$ws.Range("D40").Value = "5.250"
$ws.Range("E40").Value = "  +11.46%  "
$ws.Range("D41").Value = "1.267"
$ws.Range("E41").Value = "  +4.05%  "
$ws.Range("D42").Value = "0.02277"
$ws.Range("E42").Value = "  +12.51%  "
$ws.Range("D43").Value = "0.2051"
$ws.Range("E43").Value = "  +11.45%  "
$ws.Range("D44").Value = "0.6038"
$ws.Range("E44").Value = "  +16.04%  "
$ws.Range("D45").Value = "0.9829"
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("D46").Value = "3.833"
$ws.Range("E46").Value = "  +9.11%  "
$ws.Range("D47").Value = "13.12"
$ws.Range("E47").Value = "  +9.08%  "
$ws.Range("D48").Value = "0.5876"
$ws.Range("E48").Value = "  +14.04%  "
$ws.Range("E49").Value = "  +5.77%  "
$ws.Range("D50").Value = "1.984"
$ws.Range("E50").Value = "  +10.40%  "
$ws.Range("D51").Value = "0.07025"
$ws.Range("E51").Value = "  +9.46%  "
